# Apply updated TPM data to Sfrp1-Fzd6 sheet (adds "Resolving-Mac" target cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sfrp1"
$ws.Cells.Item(2, 3).Value = "Fzd6"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(2, 7).Value = [double]"0.2147063333333333"
$ws.Cells.Item(2, 8).Value = [double]"0.644119"
$ws.Cells.Item(2, 9).Value = [double]"0.00381773955517184"
$ws.Cells.Item(2, 10).Value = [double]"0.00381773955517184"
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = [double]"19.84402066666667"
$ws.Cells.Item(2, 14).Value = [double]"59.532062"
$ws.Cells.Item(2, 15).Value = [double]"0.8956779416773022"
$ws.Cells.Item(2, 16).Value = [double]"0.8956779416773021"
$ws.Cells.Item(2, 17).Value = [double]"4.260636915930888"
$ws.Cells.Item(2, 18).Value = [double]"38.345732243378"
$ws.Cells.Item(2, 19).Value = [double]"0.003419465106636333"
$ws.Cells.Item(2, 20).Value = [double]"0.003419465106636333"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sfrp1"
$ws.Cells.Item(3, 3).Value = "Fzd6"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(3, 7).Value = [double]"0.2147063333333333"
$ws.Cells.Item(3, 8).Value = [double]"0.644119"
$ws.Cells.Item(3, 9).Value = [double]"0.00381773955517184"
$ws.Cells.Item(3, 10).Value = [double]"0.00381773955517184"
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = [double]"1.745879666666666"
$ws.Cells.Item(3, 14).Value = [double]"5.237639"
$ws.Cells.Item(3, 15).Value = [double]"0.07880186845818919"
$ws.Cells.Item(3, 16).Value = [double]"0.07880186845818919"
$ws.Cells.Item(3, 17).Value = [double]"0.3748514216712222"
$ws.Cells.Item(3, 18).Value = [double]"3.373662795041"
$ws.Cells.Item(3, 19).Value = [double]"0.000300845010234277"
$ws.Cells.Item(3, 20).Value = [double]"0.0003008450102342771"

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sfrp1"
$ws.Cells.Item(4, 3).Value = "Fzd6"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(4, 7).Value = [double]"0.2147063333333333"
$ws.Cells.Item(4, 8).Value = [double]"0.644119"
$ws.Cells.Item(4, 9).Value = [double]"0.00381773955517184"
$ws.Cells.Item(4, 10).Value = [double]"0.00381773955517184"
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = [double]"0.5556126666666666"
$ws.Cells.Item(4, 14).Value = [double]"1.666838"
$ws.Cells.Item(4, 15).Value = [double]"0.02507808362071368"
$ws.Cells.Item(4, 16).Value = [double]"0.02507808362071367"
$ws.Cells.Item(4, 17).Value = [double]"0.1192935584135555"
$ws.Cells.Item(4, 18).Value = [double]"1.073642025722"
$ws.Cells.Item(4, 19).Value = [double]"9.574159180670563e-05"
$ws.Cells.Item(4, 20).Value = [double]"9.574159180670563e-05"

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Sfrp1"
$ws.Cells.Item(5, 3).Value = "Fzd6"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(5, 7).Value = [double]"0.2147063333333333"
$ws.Cells.Item(5, 8).Value = [double]"0.644119"
$ws.Cells.Item(5, 9).Value = [double]"0.00381773955517184"
$ws.Cells.Item(5, 10).Value = [double]"0.00381773955517184"
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5, 13).Value = [double]"0.009795"
$ws.Cells.Item(5, 14).Value = [double]"0.029385"
$ws.Cells.Item(5, 15).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(5, 16).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(5, 17).Value = [double]"0.002103048535"
$ws.Cells.Item(5, 18).Value = [double]"0.018927436815"
$ws.Cells.Item(5, 19).Value = [double]"1.68784649452439e-06"
$ws.Cells.Item(5, 20).Value = [double]"1.68784649452439e-06"

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Sfrp1"
$ws.Cells.Item(6, 3).Value = "Fzd6"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = [double]"55.43187966666667"
$ws.Cells.Item(6, 8).Value = [double]"166.295639"
$ws.Cells.Item(6, 9).Value = [double]"0.9856461909412342"
$ws.Cells.Item(6, 10).Value = [double]"0.9856461909412343"
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = [double]"19.84402066666667"
$ws.Cells.Item(6, 14).Value = [double]"59.532062"
$ws.Cells.Item(6, 15).Value = [double]"0.8956779416773022"
$ws.Cells.Item(6, 16).Value = [double]"0.8956779416773021"
$ws.Cells.Item(6, 17).Value = [double]"1099.991365697513"
$ws.Cells.Item(6, 18).Value = [double]"9899.922291277617"
$ws.Cells.Item(6, 19).Value = [double]"0.8828215515243178"
$ws.Cells.Item(6, 20).Value = [double]"0.8828215515243178"

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Sfrp1"
$ws.Cells.Item(7, 3).Value = "Fzd6"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = [double]"55.43187966666667"
$ws.Cells.Item(7, 8).Value = [double]"166.295639"
$ws.Cells.Item(7, 9).Value = [double]"0.9856461909412342"
$ws.Cells.Item(7, 10).Value = [double]"0.9856461909412343"
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = [double]"1.745879666666666"
$ws.Cells.Item(7, 14).Value = [double]"5.237639"
$ws.Cells.Item(7, 15).Value = [double]"0.07880186845818919"
$ws.Cells.Item(7, 16).Value = [double]"0.07880186845818919"
$ws.Cells.Item(7, 17).Value = [double]"96.77739159514677"
$ws.Cells.Item(7, 18).Value = [double]"870.996524356321"
$ws.Cells.Item(7, 19).Value = [double]"0.07767076148486636"
$ws.Cells.Item(7, 20).Value = [double]"0.07767076148486637"

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Sfrp1"
$ws.Cells.Item(8, 3).Value = "Fzd6"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = [double]"55.43187966666667"
$ws.Cells.Item(8, 8).Value = [double]"166.295639"
$ws.Cells.Item(8, 9).Value = [double]"0.9856461909412342"
$ws.Cells.Item(8, 10).Value = [double]"0.9856461909412343"
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = [double]"0.5556126666666666"
$ws.Cells.Item(8, 14).Value = [double]"1.666838"
$ws.Cells.Item(8, 15).Value = [double]"0.02507808362071368"
$ws.Cells.Item(8, 16).Value = [double]"0.02507808362071367"
$ws.Cells.Item(8, 17).Value = [double]"30.79865447994244"
$ws.Cells.Item(8, 18).Value = [double]"277.1878903194819"
$ws.Cells.Item(8, 19).Value = [double]"0.02471811759686219"
$ws.Cells.Item(8, 20).Value = [double]"0.02471811759686219"

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Sfrp1"
$ws.Cells.Item(9, 3).Value = "Fzd6"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = [double]"55.43187966666667"
$ws.Cells.Item(9, 8).Value = [double]"166.295639"
$ws.Cells.Item(9, 9).Value = [double]"0.9856461909412342"
$ws.Cells.Item(9, 10).Value = [double]"0.9856461909412343"
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(9, 13).Value = [double]"0.009795"
$ws.Cells.Item(9, 14).Value = [double]"0.029385"
$ws.Cells.Item(9, 15).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(9, 16).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(9, 17).Value = [double]"0.542955261335"
$ws.Cells.Item(9, 18).Value = [double]"4.886597352015"
$ws.Cells.Item(9, 19).Value = [double]"0.00043576033518782"
$ws.Cells.Item(9, 20).Value = [double]"0.00043576033518782"

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Sfrp1"
$ws.Cells.Item(10, 3).Value = "Fzd6"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = [double]"0.5925393333333333"
$ws.Cells.Item(10, 8).Value = [double]"1.777618"
$ws.Cells.Item(10, 9).Value = [double]"0.01053606950359399"
$ws.Cells.Item(10, 10).Value = [double]"0.01053606950359399"
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = [double]"19.84402066666667"
$ws.Cells.Item(10, 14).Value = [double]"59.532062"
$ws.Cells.Item(10, 15).Value = [double]"0.8956779416773022"
$ws.Cells.Item(10, 16).Value = [double]"0.8956779416773021"
$ws.Cells.Item(10, 17).Value = [double]"11.75836277647955"
$ws.Cells.Item(10, 18).Value = [double]"105.825264988316"
$ws.Cells.Item(10, 19).Value = [double]"0.009436925046348056"
$ws.Cells.Item(10, 20).Value = [double]"0.009436925046348058"

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Sfrp1"
$ws.Cells.Item(11, 3).Value = "Fzd6"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = [double]"0.5925393333333333"
$ws.Cells.Item(11, 8).Value = [double]"1.777618"
$ws.Cells.Item(11, 9).Value = [double]"0.01053606950359399"
$ws.Cells.Item(11, 10).Value = [double]"0.01053606950359399"
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = [double]"1.745879666666666"
$ws.Cells.Item(11, 14).Value = [double]"5.237639"
$ws.Cells.Item(11, 15).Value = [double]"0.07880186845818919"
$ws.Cells.Item(11, 16).Value = [double]"0.07880186845818919"
$ws.Cells.Item(11, 17).Value = [double]"1.034502373766889"
$ws.Cells.Item(11, 18).Value = [double]"9.310521363902"
$ws.Cells.Item(11, 19).Value = [double]"0.000830261963088552"
$ws.Cells.Item(11, 20).Value = [double]"0.0008302619630885522"

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Sfrp1"
$ws.Cells.Item(12, 3).Value = "Fzd6"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = [double]"0.5925393333333333"
$ws.Cells.Item(12, 8).Value = [double]"1.777618"
$ws.Cells.Item(12, 9).Value = [double]"0.01053606950359399"
$ws.Cells.Item(12, 10).Value = [double]"0.01053606950359399"
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = [double]"0.5556126666666666"
$ws.Cells.Item(12, 14).Value = [double]"1.666838"
$ws.Cells.Item(12, 15).Value = [double]"0.02507808362071368"
$ws.Cells.Item(12, 16).Value = [double]"0.02507808362071367"
$ws.Cells.Item(12, 17).Value = [double]"0.3292223590982222"
$ws.Cells.Item(12, 18).Value = [double]"2.963001231883999"
$ws.Cells.Item(12, 19).Value = [double]"0.0002642244320447812"
$ws.Cells.Item(12, 20).Value = [double]"0.0002642244320447812"

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Sfrp1"
$ws.Cells.Item(13, 3).Value = "Fzd6"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = [double]"0.5925393333333333"
$ws.Cells.Item(13, 8).Value = [double]"1.777618"
$ws.Cells.Item(13, 9).Value = [double]"0.01053606950359399"
$ws.Cells.Item(13, 10).Value = [double]"0.01053606950359399"
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = [double]"0.3333333333333333"
$ws.Cells.Item(13, 13).Value = [double]"0.009795"
$ws.Cells.Item(13, 14).Value = [double]"0.029385"
$ws.Cells.Item(13, 15).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(13, 16).Value = [double]"0.0004421062437949407"
$ws.Cells.Item(13, 17).Value = [double]"0.00580392277"
$ws.Cells.Item(13, 18).Value = [double]"0.05223530493"
$ws.Cells.Item(13, 19).Value = [double]"4.658062112596363e-06"
$ws.Cells.Item(13, 20).Value = [double]"4.658062112596364e-06"

Write-Host "done"